# Insert a new "Match ID" column at the very beginning of the sheet.
# This shifts all existing columns (A:W) one place to the right (B:X)
# and leaves a blank column A ready to be populated.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Columns.Item(1).Insert()

# Header label for the new column (row 2 holds the human-readable headers).
$ws.Range("A2").Value = "Match ID"

# Give the header cell (and the rest of the new column down through the
# last visible data row) the same bold font used by the other header
# cells, but without the border/centering used on the grouped headers.
$ws.Range("A2:A19").Font.Bold = $true

# Fill in the Match ID value for every player data row.
for ($r = 4; $r -le 20; $r++) {
    $ws.Cells.Item($r, 1).Value = 29
}

# Reflect the selection that was left active in the saved workbook.
[void]$ws.Range("A2:A19").Select()
